$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing text storage (so numeric-looking
# strings like "0.999" or "71.907.70" are not silently coerced to numbers),
# and restore the cell's original style afterwards (NumberFormat="@" bumps
# the style index otherwise).
function Set-TextCell($addr, $val) {
    $sty = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $sty
}

# Row 2: Bitcoin
Set-TextCell "D2" "71.907.70"
Set-TextCell "E2" "  -0.64%  "

# Row 3: Ethereum
Set-TextCell "D3" "3.898.11"
Set-TextCell "E3" "  -1.57%  "

# Row 4: TetherUSD
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.11%  "

# Row 5: BNB
Set-TextCell "D5" "599.37"
Set-TextCell "E5" "  +2.95%  "

# Row 6: Solana
Set-TextCell "D6" "169.83"
Set-TextCell "E6" "  +9.76%  "

# Row 7: XRP
Set-TextCell "D7" "0.680"
Set-TextCell "E7" "  +1.11%  "

# Row 8: USDC
Set-TextCell "E8" "  +0.10%  "

# Row 9: Cardano
Set-TextCell "D9" "0.764"
Set-TextCell "E9" "  +3.43%  "

# Row 10: Dogecoin
Set-TextCell "D10" "0.179"
Set-TextCell "E10" "  +7.75%  "

# Row 11: Avalanche
Set-TextCell "D11" "54.52"
Set-TextCell "E11" "  +3.44%  "

# Row 12: ShibaInu
Set-TextCell "D12" "0.0000324"
Set-TextCell "E12" "  +3.34%  "

# Row 13: Polkadot
Set-TextCell "D13" "11.36"
Set-TextCell "E13" "  +6.34%  "

# Row 14: WrappedliquidstakedEther2.0
Set-TextCell "D14" "4.510.22"
Set-TextCell "E14" "  -1.88%  "

# Row 15: WrappedEther
Set-TextCell "D15" "3.914.85"
Set-TextCell "E15" "  -1.21%  "

# Row 16: Chainlink
Set-TextCell "D16" "20.99"
Set-TextCell "E16" "  +3.60%  "

# Row 17: Uniswap
Set-TextCell "D17" "14.00"
Set-TextCell "E17" "  +1.02%  "

# Row 18: Polygon
Set-TextCell "E18" "  -2.55%  "

# Row 19: TRON
Set-TextCell "E19" "  -1.75%  "

# Row 20: WrappedBTC
Set-TextCell "D20" "71.488.88"
Set-TextCell "E20" "  -0.85%  "

# Row 21: BitcoinCash
Set-TextCell "D21" "438.42"
Set-TextCell "E21" "  +3.24%  "

# Row 22: PancakeSwap
Set-TextCell "D22" "4.75"
Set-TextCell "E22" "  +2.84%  "

# Row 23: Litecoin
Set-TextCell "D23" "95.14"
Set-TextCell "E23" "  +0.33%  "

# Row 24: ImmutableX
Set-TextCell "D24" "3.32"
Set-TextCell "E24" "  -2.51%  "

# Row 25: InternetComputer(DFINITY)
Set-TextCell "D25" "13.93"
Set-TextCell "E25" "  -1.33%  "

# Row 26: Toncoin
Set-TextCell "D26" "4.15"
Set-TextCell "E26" "  -3.40%  "

# Row 27: RenderToken
Set-TextCell "D27" "11.11"
Set-TextCell "E27" "  -0.40%  "

# Row 28: LEO
Set-TextCell "D28" "5.94"
Set-TextCell "E28" "  +0.40%  "

# Row 29: Filecoin
Set-TextCell "D29" "10.28"
Set-TextCell "E29" "  -3.58%  "

# Row 30: EthereumClassic
Set-TextCell "D30" "35.30"
Set-TextCell "E30" "  -2.02%  "

# Row 31: NEARProtocol
Set-TextCell "D31" "7.97"
Set-TextCell "E31" "  +3.46%  "

# Row 32: InjectiveProtocol
Set-TextCell "D32" "52.62"
Set-TextCell "E32" "  +8.02%  "

# Row 33: Cosmos
Set-TextCell "D33" "13.65"
Set-TextCell "E33" "  +2.55%  "

# Row 34: Hedera
Set-TextCell "E34" "  -2.36%  "

# Row 35: PEPE
Set-TextCell "D35" "0.0₂01000"
Set-TextCell "E35" "  +18.38%  "

# Row 36: OKB
Set-TextCell "D36" "68.99"
Set-TextCell "E36" "  +1.52%  "

# Row 37: Bittensor
Set-TextCell "D37" "624.07"
Set-TextCell "E37" "  -7.71%  "

# Row 38: TheGraph
Set-TextCell "D38" "0.422"
Set-TextCell "E38" "  -2.49%  "

# Row 39: Dai
Set-TextCell "D39" "0.999"
Set-TextCell "E39" "  -0.08%  "

# Row 40: ThetaToken
Set-TextCell "D40" "3.32"
Set-TextCell "E40" "  +1.21%  "

# Row 41: was FirstDigitalUSD -> now Kaspa
Set-TextCell "B41" "Kaspa"
Set-TextCell "C41" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D41" "0.143"
Set-TextCell "E41" "  -0.65%  "

# Row 42: was Kaspa -> now FirstDigitalUSD
Set-TextCell "B42" "FirstDigitalUSD"
Set-TextCell "C42" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D42" "0.998"
Set-TextCell "E42" "  -0.21%  "

# Row 43: dogwifhat
Set-TextCell "E43" "  +42.42%  "

# Row 44: VeChain
Set-TextCell "D44" "0.0472"
Set-TextCell "E44" "  -1.51%  "

# Row 45: THORChain
Set-TextCell "D45" "10.27"
Set-TextCell "E45" "  -5.30%  "

# Row 46: Fetch.AI
Set-TextCell "D46" "2.67"
Set-TextCell "E46" "  -0.53%  "

# Row 47: Stellar
Set-TextCell "E47" "  -0.97%  "

# Row 48: ApeXProtocol
Set-TextCell "D48" "3.39"
Set-TextCell "E48" "  -0.32%  "

# Row 49: was WEMIXToken -> now Maker
Set-TextCell "B49" "Maker"
Set-TextCell "C49" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell "D49" "2.858.52"
Set-TextCell "E49" "  +3.48%  "

# Row 50: was Maker -> now WEMIXToken
Set-TextCell "B50" "WEMIXToken"
Set-TextCell "C50" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell "D50" "2.78"
Set-TextCell "E50" "  -16.81%  "

# Row 51: FLOKI
Set-TextCell "D51" "0.000277"
Set-TextCell "E51" "  +3.55%  "
